$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Row 10
$tbl.ListRows.Add() | Out-Null
$ws.Range("A10").Value = 45368
$ws.Range("B10").Value = "Spez 1"
$ws.Range("C10").Value = "Kaffee und Gebäck"
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 16

# Row 11
$tbl.ListRows.Add() | Out-Null
$ws.Range("A11").Value = 45368
$ws.Range("B11").Value = "Spez 2"
$ws.Range("C11").Value = "Gebäck"
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 3

# Row 12
$tbl.ListRows.Add() | Out-Null
$ws.Range("A12").Value = 45368
$ws.Range("B12").Value = "Spez 3"
$ws.Range("C12").Value = "Kollekte"
$ws.Range("D12").Value = 240
$ws.Range("E12").Value = 1

# Excel leaves the table boundary one row beyond the last populated data row
# (matches a manual drag-resize of the table down one extra row).
$tbl.Resize($ws.Range("A1:E13")) | Out-Null

$ws.Range("E17").Select() | Out-Null
